$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update score values (Web Design: 8 -> 7, AngularJS Project Structure: 4 -> 3,
#     Authorization Checks: blank -> 5). C51 SUM(C6:C50) will recalc automatically.
$ws.Range("C11").Value = 7
$ws.Range("C12").Value = 3
$ws.Range("C32").Value = 5

# --- Update the view: scroll the window so row 40 is at the top-left, and move
#     the active selection to D55.
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("D55").Select()
